$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 146, shifting rows
# 146-159 down to 148-161 (values unchanged), so that the new rows
# 146 and 147 can hold the newly reported "Lane Late" / "Navel Late"
# weekly price records.
$ws.Rows("146:147").Insert()

# Row 146: new "Lane Late" record
$ws.Range("A146").Value = 11
$ws.Range("B146").Value = "Vega Monumental Concepción"
$ws.Range("C146").Value = "Bíobío"
$ws.Range("D146").Value = 44491
$ws.Range("E146").Value = 8
$ws.Range("F146").Value = "Fruta"
$ws.Range("G146").Value = 100102
$ws.Range("H146").Value = "Cítricos"
$ws.Range("I146").Value = 100102005
$ws.Range("J146").Value = "Naranja"
$ws.Range("K146").Value = "Lane Late"
$ws.Range("L146").Value = "Primera"
$ws.Range("M146").Value = 350
$ws.Range("N146").Value = 7500
$ws.Range("O146").Value = 8000
$ws.Range("P146").Value = 7786
$ws.Range("Q146").Value = "$/caja 15 kilos empedrada"
$ws.Range("R146").Value = "Región de O'Higgins"
$ws.Range("S146").Value = 519
$ws.Range("T146").Value = 15

# Row 147: new "Navel Late" record
$ws.Range("A147").Value = 11
$ws.Range("B147").Value = "Vega Monumental Concepción"
$ws.Range("C147").Value = "Bíobío"
$ws.Range("D147").Value = 44491
$ws.Range("E147").Value = 8
$ws.Range("F147").Value = "Fruta"
$ws.Range("G147").Value = 100102
$ws.Range("H147").Value = "Cítricos"
$ws.Range("I147").Value = 100102005
$ws.Range("J147").Value = "Naranja"
$ws.Range("K147").Value = "Navel Late"
$ws.Range("L147").Value = "Primera"
$ws.Range("M147").Value = 350
$ws.Range("N147").Value = 7500
$ws.Range("O147").Value = 8000
$ws.Range("P147").Value = 7786
$ws.Range("Q147").Value = "$/caja 15 kilos empedrada"
$ws.Range("R147").Value = "Región de O'Higgins"
$ws.Range("S147").Value = 519
$ws.Range("T147").Value = 15
